$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell, preventing Excel from
# auto-converting numeric-looking strings (e.g. "1.000", "48.00") into
# real numbers, which would lose the formatted/trailing-zero text.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$ws.Range("D2").Value = "27.292.88"
$ws.Range("E2").Value = "  -4.46%  "
$ws.Range("D3").Value = "1.857.34"
$ws.Range("E3").Value = "  -5.55%  "
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  -1.20%  "
Set-TextValue $ws.Range("D5") "321.99"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  -1.03%  "
Set-TextValue $ws.Range("D7") "0.4496"
$ws.Range("E7").Value = "  -6.29%  "
Set-TextValue $ws.Range("D8") "0.3847"
$ws.Range("E8").Value = "  -5.43%  "
Set-TextValue $ws.Range("D9") "48.00"
Set-TextValue $ws.Range("D10") "0.07873"
$ws.Range("E10").Value = "  -7.59%  "
Set-TextValue $ws.Range("D11") "1.018"
$ws.Range("E11").Value = "  -4.29%  "
Set-TextValue $ws.Range("D12") "21.32"
$ws.Range("E12").Value = "  -5.13%  "
$ws.Range("D13").Value = "1.865.58"
$ws.Range("E13").Value = "  -7.06%  "
Set-TextValue $ws.Range("D14") "7.155"
$ws.Range("E14").Value = "  -6.38%  "
Set-TextValue $ws.Range("D15") "5.871"
$ws.Range("E15").Value = "  -5.40%  "
Set-TextValue $ws.Range("D16") "1.000"
$ws.Range("E16").Value = "  -1.27%  "
Set-TextValue $ws.Range("D17") "0.00001033"
$ws.Range("E17").Value = "  -3.93%  "
Set-TextValue $ws.Range("D18") "85.39"
$ws.Range("E18").Value = "  -6.54%  "
Set-TextValue $ws.Range("D19") "0.06532"
$ws.Range("E19").Value = "  -1.82%  "
Set-TextValue $ws.Range("D20") "16.95"
$ws.Range("E20").Value = "  -9.06%  "
Set-TextValue $ws.Range("D21") "1.000"
$ws.Range("E21").Value = "  -1.09%  "
Set-TextValue $ws.Range("D22") "5.507"
$ws.Range("E22").Value = "  -6.36%  "
$ws.Range("D23").Value = "27.297.25"
$ws.Range("E23").Value = "  -4.56%  "
Set-TextValue $ws.Range("D24") "10.73"
$ws.Range("E24").Value = "  -6.86%  "
Set-TextValue $ws.Range("D25") "2.266"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").Value = "2.082.56"
$ws.Range("E26").Value = "  -7.11%  "
Set-TextValue $ws.Range("D27") "151.74"
$ws.Range("E27").Value = "  -3.01%  "
Set-TextValue $ws.Range("D28") "19.69"
$ws.Range("E28").Value = "  -3.56%  "
Set-TextValue $ws.Range("D29") "2.054"
$ws.Range("E29").Value = "  -6.01%  "
Set-TextValue $ws.Range("D30") "5.462"
$ws.Range("E30").Value = "  -7.46%  "
Set-TextValue $ws.Range("D31") "120.27"
$ws.Range("E31").Value = "  -3.76%  "
Set-TextValue $ws.Range("D32") "1.473"
$ws.Range("E32").Value = "  +0.68%  "
Set-TextValue $ws.Range("D33") "0.09277"
$ws.Range("E33").Value = "  -4.22%  "
Set-TextValue $ws.Range("D34") "0.9322"
$ws.Range("E34").Value = "  -6.03%  "
Set-TextValue $ws.Range("D35") "3.605"
$ws.Range("E35").Value = "  -2.75%  "
Set-TextValue $ws.Range("D36") "5.270"
$ws.Range("E36").Value = "  -6.69%  "
Set-TextValue $ws.Range("D37") "0.02223"
$ws.Range("E37").Value = "  -5.03%  "
Set-TextValue $ws.Range("D38") "0.05980"
$ws.Range("E38").Value = "  -4.32%  "
Set-TextValue $ws.Range("D39") "1.208"
$ws.Range("E39").Value = "  -3.96%  "
Set-TextValue $ws.Range("D40") "8.249"
$ws.Range("E40").Value = "  -9.82%  "
Set-TextValue $ws.Range("D41") "0.9998"
$ws.Range("E41").Value = "  -1.07%  "
Set-TextValue $ws.Range("D42") "0.5897"
$ws.Range("E42").Value = "  -5.66%  "
Set-TextValue $ws.Range("D43") "0.1882"
$ws.Range("E43").Value = "  -2.21%  "
Set-TextValue $ws.Range("D44") "10.11"
$ws.Range("E44").Value = "  -10.00%  "
Set-TextValue $ws.Range("D45") "1.254"
$ws.Range("E45").Value = "  -8.10%  "
Set-TextValue $ws.Range("D46") "0.5623"
$ws.Range("E46").Value = "  -5.85%  "
Set-TextValue $ws.Range("D47") "11.95"
$ws.Range("E47").Value = "  -9.16%  "
Set-TextValue $ws.Range("D48") "3.354"
$ws.Range("E48").Value = "  -1.84%  "
Set-TextValue $ws.Range("D49") "1.918"
$ws.Range("E49").Value = "  -7.45%  "
$ws.Range("E50").Value = "  -0.50%  "
Set-TextValue $ws.Range("D51") "108.01"
$ws.Range("E51").Value = "  -3.20%  "
